# #5: fund, bonds, otherbonds, antique done
#
# The "基金受益憑證" (fund) worksheet only had the "brief" columns
# (name/owner/dealer/quantity/face_value/currency/total). Bring it in line
# with the other sheets in this workbook by adding the shared metadata
# columns (property_category, category, date, legislator_name,
# legislator_id, source_file, index) and backfilling the header row, which
# was previously missing its own "name"/"owner"/"dealer" labels (it was
# reusing row-2's data values by mistake).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# ---- header row (row 1) ----
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "dealer"
$ws.Range("E1").Value = "quantity"
$ws.Range("F1").Value = "face_value"
$ws.Range("G1").Value = "currency"
$ws.Range("H1").Value = "total"
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

# header row style (bold + border, matching B1:H1) needs to cover the new
# I1:O1 cells too
$ws.Range("B1").Copy()
$ws.Range("I1:O1").PasteSpecial(-4122)

# ---- row 2 ----
$ws.Range("B2").Value = "台灣工銀大眾基金"
$ws.Range("D2").Value = "台灣工銀證券投信公司"
$ws.Range("I2").Value = "fund"
$ws.Range("J2").Value = "normal"
$ws.Range("L2").Value = "黃昭順"
$ws.Range("M2").Value = 665
$ws.Range("N2").Value = "tmp4c4f1"
$ws.Range("O2").Value = 91

# ---- row 3 ----
$ws.Range("B3").Value = "國泰全球貨幣"
$ws.Range("D3").Value = "國泰證券投信公司"
$ws.Range("I3").Value = "fund"
$ws.Range("J3").Value = "normal"
$ws.Range("L3").Value = "黃昭順"
$ws.Range("M3").Value = 665
$ws.Range("N3").Value = "tmp4c4f1"
$ws.Range("O3").Value = 92

# "2012-04-26" needs to land in K2/K3 as plain text (matching every other
# sheet's "date" column), not get auto-converted into a date serial number
# by a straight .Value assignment -- pull it in as a value-only paste from
# a cell that already holds it as text (sheet "汽車", J2).
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("J2").Copy()
$ws.Range("K2").PasteSpecial(-4163)
$ws1.Range("J2").Copy()
$ws.Range("K3").PasteSpecial(-4163)

# data-row style (plain, matching B2:H3) needs to cover the new I2:O3 cells
$ws.Range("B2").Copy()
$ws.Range("I2:O3").PasteSpecial(-4122)
